$wb = $excel.ActiveWorkbook

# --- Sheet "About": insert two new explanatory lines + one blank row
# after the existing "Currency Year" section header (old row 17), pushing
# the old rows 18-25 down to 21-28.
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Rows.Item(18).Insert()
$wsAbout.Rows.Item(18).Insert()
$wsAbout.Rows.Item(18).Insert()

# The inserted rows inherit the bold "section header" formatting of row 17
# above them; these new lines are plain body text like the paragraph below,
# so drop the inherited bold weight.
$wsAbout.Range("A18:A20").Font.Bold = $false

$wsAbout.Range("A18").Value = "The model uses LDVs elasticity for all vehicle types because no data on price elasticity"
$wsAbout.Range("A19").Value = "of other vehicle types with respect to fuel economy is available."

# --- Sheet "Calculations": the formula in B5 pointed at About!A24 (the
# 2010->2012 dollar conversion factor); after the row insertion above that
# value now lives three rows further down, at About!A27.
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsCalc.Range("B5").Formula = "=B4*About!A27"

# --- Sheet "EoVPwFE": clarify the elasticity label and wrap/right-align it.
$wsEoV = $wb.Worksheets.Item("EoVPwFE")
$wsEoV.Range("B1").Value = "Elasticity (dimensionless)"
$wsEoV.Range("B1").HorizontalAlignment = -4152
$wsEoV.Range("B1").WrapText = $true
$wsEoV.Range("B1").EntireRow.RowHeight = 45

# Restore the on-screen selections to match the saved workbook view state.
$wsEoV.Range("B1").Select()

$wsAbout.Activate()
$wsAbout.Range("A20:XFD23").Select()
